$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 926825.75
$ws.Range("I11").Value = 926825.75
$ws.Range("K11").Value = 926825.75
$ws.Range("M11").Value = -926685.75
$ws.Range("H109").Value = 85791.664
$ws.Range("J109").Value = 85791.664
$ws.Range("L109").Value = 85791.664
$ws.Range("N109").Value = -88565.664
$ws.Range("H116").Value = 1673624.6
$ws.Range("I116").Value = 7188.9
$ws.Range("K116").Value = 7188.9
$ws.Range("M116").Value = -3746.9
$ws.Range("H125").Value = 6405.7
$ws.Range("I125").Value = 4817.4
$ws.Range("J125").Value = 7994
$ws.Range("K125").Value = 43356.6
$ws.Range("L125").Value = 71946
$ws.Range("M125").Value = -40896.6
$ws.Range("N125").Value = -76866
$ws.Range("H133").Value = 89838.75
$ws.Range("J133").Value = 89838.75
$ws.Range("L133").Value = 89838.75
$ws.Range("N133").Value = -99958.75
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
$ws.Range("H140").Value = 91491.664
$ws.Range("J140").Value = 91990
$ws.Range("L140").Value = 91990
$ws.Range("N140").Value = -102350
$ws.Range("H141").Value = 2807.68
$ws.Range("I141").Value = 2828.0908
$ws.Range("J141").Value = 2658
$ws.Range("K141").Value = 8484.2724
$ws.Range("L141").Value = 7974
$ws.Range("M141").Value = -3304.2724
$ws.Range("N141").Value = -18334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2622.8333
$ws.Range("I61").Value = 2119.4443
$ws.Range("K61").Value = 2119.4443
$ws.Range("M61").Value = -1907.4443
$ws.Range("H74").Value = 2500.4443
$ws.Range("I74").Value = 3685.6667
$ws.Range("J74").Value = 2263.4
$ws.Range("K74").Value = 3685.6667
$ws.Range("L74").Value = 2263.4
$ws.Range("M74").Value = -2811.6667
$ws.Range("N74").Value = -4011.4
$ws.Range("H77").Value = 2500.4443
$ws.Range("I77").Value = 3685.6667
$ws.Range("J77").Value = 2263.4
$ws.Range("K77").Value = 18428.3335
$ws.Range("L77").Value = 11317
$ws.Range("M77").Value = -14060.3335
$ws.Range("N77").Value = -20053
$ws.Range("H122").Value = 2506.0303
$ws.Range("I122").Value = 2478.926
$ws.Range("J122").Value = 2628
$ws.Range("K122").Value = 7436.778
$ws.Range("L122").Value = 7884
$ws.Range("M122").Value = -4986.778
$ws.Range("N122").Value = -12784
$ws.Range("H132").Value = 1808.0385
$ws.Range("J132").Value = 4000.3333
$ws.Range("L132").Value = 12000.9999
$ws.Range("N132").Value = -17060.9999
$ws.Range("H136").Value = 2622.8333
$ws.Range("I136").Value = 2119.4443
$ws.Range("K136").Value = 6358.3329
$ws.Range("M136").Value = -3808.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 33096.668
$ws.Range("J2").Value = 33096.668
$ws.Range("L2").Value = 33096.668
$ws.Range("N2").Value = -33322.668
$ws.Range("H20").Value = 54718.105
$ws.Range("I20").Value = 78666.53999999999
$ws.Range("J20").Value = 2829.8333
$ws.Range("K20").Value = 78666.53999999999
$ws.Range("L20").Value = 2829.8333
$ws.Range("M20").Value = -78419.53999999999
$ws.Range("N20").Value = -3323.8333
$ws.Range("H88").Value = 17561
$ws.Range("J88").Value = 17561
$ws.Range("L88").Value = 17561
$ws.Range("N88").Value = -18373
$ws.Range("H91").Value = 17561
$ws.Range("J91").Value = 17561
$ws.Range("L91").Value = 17561
$ws.Range("N91").Value = -20369
$ws.Range("H114").Value = 89095.57000000001
$ws.Range("J114").Value = 88674.664
$ws.Range("L114").Value = 88674.664
$ws.Range("N114").Value = -97352.664
$ws.Range("H116").Value = 73774
$ws.Range("J116").Value = 73774
$ws.Range("L116").Value = 73774
$ws.Range("N116").Value = -82952
$ws.Range("H118").Value = 77626.78
$ws.Range("J118").Value = 74886.75
$ws.Range("L118").Value = 74886.75
$ws.Range("N118").Value = -78200.75
$ws.Range("H134").Value = 1442.0426
$ws.Range("I134").Value = 1111.6097
$ws.Range("K134").Value = 3334.8291
$ws.Range("M134").Value = -799.8290999999999
$ws.Range("H135").Value = 47602.8
$ws.Range("J135").Value = 47602.8
$ws.Range("L135").Value = 47602.8
$ws.Range("N135").Value = -57742.8
$ws.Range("H138").Value = 96368
$ws.Range("J138").Value = 96368
$ws.Range("L138").Value = 96368
$ws.Range("N138").Value = -106648

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1878.5555
$ws.Range("I16").Value = 1604.3334
$ws.Range("J16").Value = 3249.6667
$ws.Range("K16").Value = 1604.3334
$ws.Range("L16").Value = 3249.6667
$ws.Range("M16").Value = -1317.3334
$ws.Range("N16").Value = -3823.6667
$ws.Range("H31").Value = 13155.134
$ws.Range("I31").Value = 1881.1578
$ws.Range("K31").Value = 1881.1578
$ws.Range("M31").Value = -1586.1578
$ws.Range("H34").Value = 13155.134
$ws.Range("I34").Value = 1881.1578
$ws.Range("K34").Value = 1881.1578
$ws.Range("M34").Value = -1679.1578
$ws.Range("H99").Value = 1647264.4
$ws.Range("I99").Value = 2674.8
$ws.Range("K99").Value = 2674.8
$ws.Range("M99").Value = -1176.8
$ws.Range("H113").Value = 1878.5555
$ws.Range("I113").Value = 1604.3334
$ws.Range("J113").Value = 3249.6667
$ws.Range("K113").Value = 1604.3334
$ws.Range("L113").Value = 3249.6667
$ws.Range("M113").Value = 565.6666
$ws.Range("N113").Value = -7589.6667
$ws.Range("H116").Value = 51678.332
$ws.Range("J116").Value = 51678.332
$ws.Range("L116").Value = 51678.332
$ws.Range("N116").Value = -60856.332
$ws.Range("H117").Value = 43603.5
$ws.Range("J117").Value = 43603.5
$ws.Range("L117").Value = 43603.5
$ws.Range("N117").Value = -52781.5
$ws.Range("H119").Value = 90552.28999999999
$ws.Range("J119").Value = 90552.28999999999
$ws.Range("L119").Value = 90552.28999999999
$ws.Range("N119").Value = -100228.29
$ws.Range("H126").Value = 1647264.4
$ws.Range("I126").Value = 2674.8
$ws.Range("K126").Value = 8024.400000000001
$ws.Range("M126").Value = -5554.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1313.5714
$ws.Range("I14").Value = 1313.5714
$ws.Range("K14").Value = 3940.7142
$ws.Range("M14").Value = -3767.7142
$ws.Range("H74").Value = 5413.75
$ws.Range("I74").Value = 2218.3333
$ws.Range("K74").Value = 6654.999899999999
$ws.Range("M74").Value = -5593.999899999999
$ws.Range("H75").Value = 82.5
$ws.Range("I75").Value = 100
$ws.Range("J75").Value = 65
$ws.Range("K75").Value = 300
$ws.Range("L75").Value = 195
$ws.Range("M75").Value = 698
$ws.Range("N75").Value = -2191
$ws.Range("H77").Value = 5413.75
$ws.Range("I77").Value = 2218.3333
$ws.Range("K77").Value = 19964.9997
$ws.Range("M77").Value = -14660.9997
$ws.Range("H78").Value = 82.5
$ws.Range("I78").Value = 100
$ws.Range("J78").Value = 65
$ws.Range("K78").Value = 900
$ws.Range("L78").Value = 585
$ws.Range("M78").Value = 4092
$ws.Range("N78").Value = -10569
$ws.Range("H121").Value = 1724.5588
$ws.Range("J121").Value = 1781.1072
$ws.Range("L121").Value = 5343.321599999999
$ws.Range("N121").Value = -7963.321599999999
$ws.Range("H132").Value = 9999
$ws.Range("I132").Value = 9999
$ws.Range("K132").Value = 89991
$ws.Range("M132").Value = -87461

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 824.5925999999999
$ws.Range("I102").Value = 627.7917
$ws.Range("K102").Value = 627.7917
$ws.Range("M102").Value = 994.2083
$ws.Range("H108").Value = 64238.184
$ws.Range("J108").Value = 64238.184
$ws.Range("L108").Value = 64238.184
$ws.Range("N108").Value = -71918.18400000001
$ws.Range("H113").Value = 2779488.8
$ws.Range("J113").Value = 5557891
$ws.Range("L113").Value = 5557891
$ws.Range("N113").Value = -5562231
$ws.Range("H122").Value = 2719.6
$ws.Range("I122").Value = 2051
$ws.Range("J122").Value = 4558.25
$ws.Range("K122").Value = 6153
$ws.Range("L122").Value = 13674.75
$ws.Range("M122").Value = -3703
$ws.Range("N122").Value = -18574.75
$ws.Range("H126").Value = 3474.3157
$ws.Range("I126").Value = 2904.5715
$ws.Range("J126").Value = 3806.6667
$ws.Range("K126").Value = 8713.7145
$ws.Range("L126").Value = 11420.0001
$ws.Range("M126").Value = -6243.7145
$ws.Range("N126").Value = -16360.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20771.52
$ws.Range("I7").Value = 10795.131
$ws.Range("K7").Value = 10795.131
$ws.Range("M7").Value = -10683.131
$ws.Range("H126").Value = 20771.52
$ws.Range("I126").Value = 10795.131
$ws.Range("K126").Value = 32385.393
$ws.Range("M126").Value = -29915.393
$ws.Range("H136").Value = 5367.609
$ws.Range("J136").Value = 4191.077
$ws.Range("L136").Value = 12573.231
$ws.Range("N136").Value = -17673.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 10205414
$ws.Range("I100").Value = 35714670
$ws.Range("K100").Value = 71429340
$ws.Range("M100").Value = -71428799
$ws.Range("H121").Value = 89812
$ws.Range("J121").Value = 89812
$ws.Range("L121").Value = 89812
$ws.Range("N121").Value = -93306
$ws.Range("H122").Value = 2967.9375
$ws.Range("I122").Value = 2219.4
$ws.Range("K122").Value = 6658.200000000001
$ws.Range("M122").Value = -4208.200000000001
